# Adding Lab 19 HumMod data, it cannot be completed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the explanatory note text first (matches shared-string insertion order) ---
$ws.Range("A11").Value = "Cannot be completed in HumMod. No way to turn off ADH formation"

# --- Label the existing (left) table as "QCP" and add a "HumMod" title for the new table ---
$ws.Range("C1").Value = "QCP"
$ws.Range("G1").Value = "Diabetes Insipidus"
$ws.Range("I1").Value = "HumMod"

# --- Copy the formatting of the original table (A2:E9) onto the new HumMod table (G2:K9) ---
$ws.Range("A2:E9").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the HumMod table's header / row-label values (mirrors the QCP table) ---
$ws.Range("G2").Value = "Time"
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1

$ws.Range("H3").Value = "Min"
$ws.Range("I3").Value = "Min"
$ws.Range("J3").Value = "Day"
$ws.Range("K3").Value = "Week"

$ws.Range("G4").Value = "Plasma [ADH](pg/mL)"
$ws.Range("G5").Value = "Plasma [Osm](mOsm/L)"
$ws.Range("G6").Value = "H2O Intake(mL/min)"
$ws.Range("G7").Value = "ECFV(L)"
$ws.Range("G8").Value = "Urine Flow(mL/min)"
$ws.Range("G9").Value = "Urine [Osm](mOsm/L)"

# H4:K9 intentionally left blank - HumMod cannot produce this data

# --- Merge the explanatory note in row 11 across A11:G11 with a yellow highlight ---
$ws.Range("A11:G11").Merge() | Out-Null

# B11:G11 (the non-anchor merged cells) get a plain yellow centered style
$ws.Range("B11:G11").Interior.Color = 65535        # yellow (RGB FFFF00)
$ws.Range("B11:G11").HorizontalAlignment = -4108   # xlCenter

# A11 gets the table's Arial 12 font (borrowed from A4, then border cleared) plus the yellow
# highlight, centered / top-aligned / wrapped text
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A11").Borders.LineStyle = -4142         # xlLineStyleNone
$ws.Range("A11").Interior.Color = 65535
$ws.Range("A11").HorizontalAlignment = -4108

# --- Merge the row-header cell for the new table, matching the QCP table's A2:A3 ---
$ws.Range("G2:G3").Merge() | Out-Null

# --- Match the final selection state ---
$ws.Range("K4").Select() | Out-Null
